# Sprint Cadence Diagrams.pptx - apply the authored changes:
#  1. Refresh the auto-date placeholder ("datetimeFigureOut" field) on every
#     slide layout and on the slide master from 8/26/2020 -> 8/27/2020.
#  2. Append a new blank slide (slide 4), using the same "Blank" layout that
#     the existing slides already use.

$p = $ppt.ActivePresentation

# --- 1. Update every "Date Placeholder" shape (slide layouts + master) ---
$master = $p.SlideMaster
$layouts = $master.CustomLayouts

for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $sh = $layout.Shapes.Item($si)
        if ($sh.Name -like "Date Placeholder*") {
            $sh.TextFrame.TextRange.Text = "8/27/2020"
        }
    }
}

for ($si = 1; $si -le $master.Shapes.Count; $si++) {
    $sh = $master.Shapes.Item($si)
    if ($sh.Name -like "Date Placeholder*") {
        $sh.TextFrame.TextRange.Text = "8/27/2020"
    }
}

# --- 2. Add a new blank slide at position 4 (same "Blank" layout as the
#         other slides in this deck) ---
$blankLayout = $layouts.Item(7)
$newSlide = $p.Slides.AddSlide(4, $blankLayout)

Write-Output "Slides: $($p.Slides.Count)"
